$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A3").Value = "2022-05-30 T 21:36:43 UTC"
$ws.Range("B2:B3").Value = 1919.638797075
$ws.Range("C2:C3").Value = 0.800712
$ws.Range("D2:D3").Value = 1.232407
